# Juno: check in to OLPRODLOC.
# Update the Chai Tea market-trends header row:
#   - Capitalize "chai" -> "Chai" in several header labels
#   - Rename "Ventas de chai preparado (unidades)" -> "Ventas predefinidas de Chai (unidades)"
#   - Rename "Involucración en redes sociales (visualizaciones)" -> "Interacción de redes sociales (visualizaciones)"
#   - Make the whole header row bold

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Ventas totales de Chai (unidades)"
$ws.Range("C1").Value = "Ventas de Chai artesanal (unidades)"
$ws.Range("D1").Value = "Ventas predefinidas de Chai (unidades)"
$ws.Range("E1").Value = "Interacción de redes sociales (visualizaciones)"
$ws.Range("F1").Value = "Búsquedas en línea de Chai"

# Make the header row text bold (A1:F1)
$ws.Range("A1:F1").Font.Bold = $true

$wb.Save()
